# Correct the "District" column (G) values to the official name
# "Uttara Kannada (Karwar)" for all data rows except row 16, which
# already contains a distinct value ("Ankola Uttara Kannada") that is
# left untouched per the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 3; $row -le 41; $row++) {
    if ($row -eq 16) {
        continue
    }
    $ws.Cells.Item($row, 7).Value = "Uttara Kannada (Karwar)"
}
